# CPMV parameters when mismatch
# Adds three new columns (T: VQ_best, U: comp_offs_y, V: comp_offs_x) and
# rotates the existing M/N/O (mv2_v / mv2_h / ... ) candidate values across
# rows 2-4, filling in the new compensation-offset data for rows 2-4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths for the new T/U/V columns ---
$ws.Columns.Item(20).ColumnWidth = 10.166666666666666
$ws.Columns.Item(21).ColumnWidth = 15.333333333333334
$ws.Columns.Item(22).ColumnWidth = 13.0

# --- New header row (row 1) ---
$ws.Range("T1").Value = "VQ_best"
$ws.Range("U1").Value = "comp_offs_y"
$ws.Range("V1").Value = "comp_offs_x"

# --- Row 2: rotate M/N/O values and add new T/U/V data ---
$ws.Range("M2").Value = 96
$ws.Range("N2").Value = -32
$ws.Range("O2").Value = 32
$ws.Range("T2").Value = 3
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0

# --- Row 3: rotate M/N values (O3 unchanged) and add new U/V data ---
$ws.Range("M3").Value = 88
$ws.Range("N3").Value = -56
$ws.Range("U3").Value = 99
$ws.Range("V3").Value = 187

# --- Row 4: rotate M/N/O values and add new U/V data ---
$ws.Range("M4").Value = 128
$ws.Range("N4").Value = -80
$ws.Range("O4").Value = 48
$ws.Range("U4").Value = 98
$ws.Range("V4").Value = 186

# --- Update the view's active selection to the newly added V3 cell ---
$ws.Range("V3").Select()
